# Updated symbol list on Fri Dec 30 15:54:41 UTC 2022 with GitHub Actions
#
# The crypto list shifted by one row (a new "LEO" entry was inserted at the
# top of the rankings, row 4, pushing HuobiToken..TigerCash down by one
# row), and the Price (column D) and a handful of other Price-only cells
# were refreshed with newer quotes. Column D holds numeric-looking values
# that are stored as plain text in this sheet, so every Price cell is
# written with a leading apostrophe to force Excel to keep it as text
# instead of silently re-typing it as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# BNB / OKB: price refresh only
Set-TextCell "D2" "245.05"
Set-TextCell "D3" "25.00"

# Row 4: HuobiToken -> LEO
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D4" "3.503"
$ws.Range("E4").Value = "3LEOLEO"

# Row 5: Cronos -> HuobiToken
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D5" "5.118"
$ws.Range("E5").Value = "4HuobiTokenHT"

# Row 6: KuCoinToken -> Cronos
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D6" "0.05656"
$ws.Range("E6").Value = "5CronosCRO"

# Row 7: GateToken -> KuCoinToken
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell "D7" "6.520"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8: MXToken -> GateToken
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D8" "2.955"
$ws.Range("E8").Value = "7GateTokenGT"

# Row 9: FTXToken -> MXToken
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D9" "0.8119"
$ws.Range("E9").Value = "8MXTokenMX"

# Row 10: WazirX -> FTXToken
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D10" "0.8361"
$ws.Range("E10").Value = "9FTXTokenFTT"

# Row 11: MandalaExchangeToken -> WazirX
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D11" "0.1333"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12: BitrueCoin -> MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.06949"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13: BitMartToken -> BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.02837"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14: BitForexToken -> BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09409"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15: One -> BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001507"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16: TigerCash -> One
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D16" "0.0005962"
$ws.Range("E16").Value = "15OneONE"

# Row 17: LEO -> TigerCash
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D17" "0.006106"
$ws.Range("E17").Value = "16TigerCashTCH"

# Remaining Price-only refreshes further down the table (no reordering).
Set-TextCell "D19" "0.3165"
Set-TextCell "D20" "0.03199"
Set-TextCell "D22" "3.743"
Set-TextCell "D23" "0.04676"
Set-TextCell "D24" "0.1357"
Set-TextCell "D26" "0.004263"
Set-TextCell "D27" "0.00009688"
Set-TextCell "D40" "0.03628"
Set-TextCell "D41" "0.006242"
Set-TextCell "D44" "0.007372"
Set-TextCell "D45" "0.00005267"
Set-TextCell "D47" "0.2197"
